# Update the existing sheet "Лист1": add price column (B) and force A3 to text
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Лист1")

# Ensure A3 is stored as text "267515" rather than a number, with no lingering
# custom number-format style left behind on the cell.
$ws1.Range("A3").NumberFormat = "@"
$ws1.Range("A3").Value = "267515"
$ws1.Range("A3").ClearFormats()

$ws1.Range("B1").Value = "87 900 ₽"
$ws1.Range("B2").Value = "Не найдено"
$ws1.Range("B3").Value = "105 678 ₽"

# Add a new worksheet "BonpetData" (placed after "Лист1") with the same data
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "BonpetData"

$ws2.Range("A1").Value = "АВР-Б-100-2-1"
$ws2.Range("B1").Value = "87 900 ₽"
$ws2.Range("A2").Value = "bababab"
$ws2.Range("B2").Value = "Не найдено"
$ws2.Range("A3").NumberFormat = "@"
$ws2.Range("A3").Value = "267515"
$ws2.Range("A3").ClearFormats()
$ws2.Range("B3").Value = "105 678 ₽"

# Keep "Лист1" as the active/selected sheet (matches original selection on A3)
$ws1.Activate() | Out-Null
$ws1.Range("A3").Select() | Out-Null
